# Updates the cryptos list worksheet to reflect the latest scraped prices
# and 1h volume percentages, matching the upstream GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assigns a value to a cell while forcing it to remain a text
# (string) cell even when the text looks like a numeric literal
# (e.g. "123.75"). We build it as a formula that yields the literal
# string, then convert the formula to its resulting value in place via
# PasteSpecial, which keeps the cell a plain string cell (no numeric
# coercion, no NumberFormat / style changes).
function Set-TextValue($cellRef, $value) {
    $escaped = $value -replace '"', '""'
    $cell = $ws.Range($cellRef)
    $cell.Formula = "=""$escaped"""
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# Plain value assignments (safe as text already: B/C strings, E percentages, D values that are not pure numeric literals)
$ws.Range("D2").Value = '43.602.23'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '2.275.53'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  -0.37%  '
$ws.Range("E5").Value = '  +6.33%  '
$ws.Range("E6").Value = '  -0.91%  '
$ws.Range("E7").Value = '  +2.24%  '
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("E9").Value = '  +1.13%  '
$ws.Range("E10").Value = '  -2.33%  '
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("E12").Value = '  +3.95%  '
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("E14").Value = '  -2.80%  '
$ws.Range("E15").Value = '  +3.10%  '
$ws.Range("D16").Value = '2.617.67'
$ws.Range("E16").Value = '  -0.80%  '
$ws.Range("D17").Value = '2.267.89'
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = '43.494.05'
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("E19").Value = '  +0.50%  '
$ws.Range("E20").Value = '  -0.84%  '
$ws.Range("E21").Value = '  -0.36%  '
$ws.Range("E22").Value = '  -0.40%  '
$ws.Range("E23").Value = '  +0.81%  '
$ws.Range("E24").Value = '  -2.03%  '
$ws.Range("E25").Value = '  -4.94%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("E26").Value = '  +2.25%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("E27").Value = '  +1.61%  '
$ws.Range("E28").Value = '  +0.76%  '
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("E32").Value = '  +0.27%  '
$ws.Range("E33").Value = '  -2.11%  '
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("E35").Value = '  +1.50%  '
$ws.Range("E36").Value = '  +11.81%  '
$ws.Range("E37").Value = '  +4.63%  '
$ws.Range("E38").Value = '  -1.87%  '
$ws.Range("E39").Value = '  -2.52%  '
$ws.Range("E40").Value = '  +4.67%  '
$ws.Range("E41").Value = '  -3.93%  '
$ws.Range("E42").Value = '  -0.31%  '
$ws.Range("E43").Value = '  -1.09%  '
$ws.Range("E44").Value = '  -0.33%  '
$ws.Range("E45").Value = '  -1.18%  '
$ws.Range("E46").Value = '  -11.26%  '
$ws.Range("E47").Value = '  +38.14%  '
$ws.Range("E48").Value = '  -2.15%  '
$ws.Range("E49").Value = '  -0.60%  '
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("E51").Value = '  -1.45%  '

# Text-safe value assignments for D-column numeric-looking values (must remain text cells)
Set-TextValue "D5" '123.75'
Set-TextValue "D6" '266.72'
Set-TextValue "D7" '0.642'
Set-TextValue "D9" '0.626'
Set-TextValue "D10" '48.25'
Set-TextValue "D11" '0.0948'
Set-TextValue "D12" '9.28'
Set-TextValue "D13" '0.107'
Set-TextValue "D14" '15.46'
Set-TextValue "D15" '0.910'
Set-TextValue "D20" '6.96'
Set-TextValue "D21" '72.23'
Set-TextValue "D22" '2.43'
Set-TextValue "D23" '235.59'
Set-TextValue "D24" '2.90'
Set-TextValue "D25" '9.49'
Set-TextValue "D26" '11.95'
Set-TextValue "D27" '1.02'
Set-TextValue "D28" '42.49'
Set-TextValue "D31" '172.55'
Set-TextValue "D32" '21.71'
Set-TextValue "D33" '0.0917'
Set-TextValue "D34" '5.74'
Set-TextValue "D38" '4.64'
Set-TextValue "D40" '2.56'
Set-TextValue "D41" '14.08'
Set-TextValue "D42" '74.04'
Set-TextValue "D43" '0.239'
Set-TextValue "D44" '0.999'
Set-TextValue "D46" '5.67'
Set-TextValue "D47" '74.66'
Set-TextValue "D48" '8.57'
Set-TextValue "D51" '101.78'
